# "aggiornamento fino a 9 agosto 2021"
# Appends 15 new daily-data rows (2021-07-26 .. 2021-08-09) to the bottom
# of the sheet, right after the last existing row (328), extending the
# used range from A1:D328 to A1:D343.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A carries a specific date/number style (bold, centered, bordered,
# custom date numFmt) on every existing data row. Copy that formatting
# from the last populated row (A328) down onto the new A329:A343 cells so
# they reuse the very same style index instead of Excel minting a new one.
$ws.Range("A328").Copy()
$ws.Range("A329:A343").PasteSpecial(-4122)

# date-serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$data = @(
  @(44403, 3, 4, 15.6561900661474),
  @(44404, 0, 4, 15.6561900661474),
  @(44405, 0, 4, 15.6561900661474),
  @(44406, 6, 10, 39.14047516536851),
  @(44407, 6, 16, 62.62476026458961),
  @(44408, 1, 17, 66.53880778112647),
  @(44409, 8, 24, 93.93714039688442),
  @(44410, 7, 28, 109.5933304630318),
  @(44411, 6, 34, 133.0776155622529),
  @(44412, 2, 36, 140.9057105953266),
  @(44413, 7, 37, 144.8197581118635),
  @(44414, 11, 42, 164.3899956945477),
  @(44415, 2, 43, 168.3040432110846),
  @(44416, 4, 39, 152.6478531449372),
  @(44417, 6, 38, 148.7338056284003)
)

$r = 329
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $r = $r + 1
}
